$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.63"
$ws.Range("E2").Value = "'0.15%"
$ws.Range("G2").Value = "'19"
$ws.Range("D3").Value = "'31.16"
$ws.Range("E3").Value = "'0.92%"
$ws.Range("G3").Value = "'19"
$ws.Range("D4").Value = "'4.959"
$ws.Range("E4").Value = "'1.34%"
$ws.Range("G4").Value = "'19"
$ws.Range("D5").Value = "'0.07477"
$ws.Range("E5").Value = "'2.76%"
$ws.Range("G5").Value = "'19"
$ws.Range("D6").Value = "'2.310"
$ws.Range("E6").Value = "'1.83%"
$ws.Range("G6").Value = "'19"
$ws.Range("D7").Value = "'7.779"
$ws.Range("E7").Value = "'1.28%"
$ws.Range("G7").Value = "'19"
$ws.Range("D8").Value = "'0.9193"
$ws.Range("E8").Value = "'2.33%"
$ws.Range("G8").Value = "'19"
$ws.Range("D9").Value = "'0.09442"
$ws.Range("E9").Value = "'19.49%"
$ws.Range("G9").Value = "'19"
$ws.Range("D10").Value = "'0.1735"
$ws.Range("E10").Value = "'4.04%"
$ws.Range("G10").Value = "'19"
$ws.Range("D11").Value = "'0.08376"
$ws.Range("E11").Value = "'4.13%"
$ws.Range("G11").Value = "'19"
$ws.Range("D12").Value = "'0.03288"
$ws.Range("E12").Value = "'6.11%"
$ws.Range("G12").Value = "'19"
$ws.Range("D13").Value = "'0.09941"
$ws.Range("E13").Value = "'-1.03%"
$ws.Range("G13").Value = "'19"
$ws.Range("D14").Value = "'0.001498"
$ws.Range("E14").Value = "'0.06%"
$ws.Range("G14").Value = "'19"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04537"
$ws.Range("E15").Value = "'0.62%"
$ws.Range("G15").Value = "'19"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005684"
$ws.Range("E16").Value = "'-1.56%"
$ws.Range("G16").Value = "'19"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.473"
$ws.Range("E17").Value = "'-0.28%"
$ws.Range("G17").Value = "'19"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.775"
$ws.Range("E18").Value = "'1.73%"
$ws.Range("G18").Value = "'19"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.195"
$ws.Range("E19").Value = "'5.63%"
$ws.Range("G19").Value = "'19"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3333"
$ws.Range("E20").Value = "'0.29%"
$ws.Range("G20").Value = "'19"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D21").Value = "'0.1309"
$ws.Range("E21").Value = "'0.76%"
$ws.Range("G21").Value = "'19"
$ws.Range("B22").Value = "MCDex"
$ws.Range("C22").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D22").Value = "'4.094"
$ws.Range("E22").Value = "'1.44%"
$ws.Range("G22").Value = "'19"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "'0.2123"
$ws.Range("E23").Value = "'1.19%"
$ws.Range("G23").Value = "'19"
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'0.88%"
$ws.Range("G24").Value = "'19"
$ws.Range("D25").Value = "'0.004304"
$ws.Range("E25").Value = "'-7.59%"
$ws.Range("G25").Value = "'19"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("G26").Value = "'19"
$ws.Range("D27").Value = "'0.0003391"
$ws.Range("E27").Value = "'0.02%"
$ws.Range("G27").Value = "'19"
$ws.Range("G28").Value = "'19"
$ws.Range("G29").Value = "'19"
$ws.Range("G30").Value = "'19"
$ws.Range("G31").Value = "'19"
$ws.Range("G32").Value = "'19"
$ws.Range("G33").Value = "'19"
$ws.Range("G34").Value = "'19"
$ws.Range("G35").Value = "'19"
$ws.Range("G36").Value = "'19"
$ws.Range("G37").Value = "'19"
$ws.Range("G38").Value = "'19"
$ws.Range("D39").Value = "'0.01621"
$ws.Range("E39").Value = "'2.41%"
$ws.Range("G39").Value = "'19"
$ws.Range("D40").Value = "'0.04579"
$ws.Range("E40").Value = "'4.24%"
$ws.Range("G40").Value = "'19"
$ws.Range("D41").Value = "'0.007505"
$ws.Range("E41").Value = "'2.94%"
$ws.Range("G41").Value = "'19"
$ws.Range("D42").Value = "'0.009831"
$ws.Range("E42").Value = "'0.78%"
$ws.Range("G42").Value = "'19"
$ws.Range("D43").Value = "'0.1360"
$ws.Range("E43").Value = "'3.51%"
$ws.Range("G43").Value = "'19"
$ws.Range("D44").Value = "'0.002218"
$ws.Range("E44").Value = "'7.28%"
$ws.Range("G44").Value = "'19"
$ws.Range("D45").Value = "'0.009040"
$ws.Range("E45").Value = "'-3.30%"
$ws.Range("G45").Value = "'19"
$ws.Range("D46").Value = "'0.00006096"
$ws.Range("E46").Value = "'6.41%"
$ws.Range("G46").Value = "'19"
$ws.Range("E47").Value = "'0.03%"
$ws.Range("G47").Value = "'19"
$ws.Range("D48").Value = "'2.551"
$ws.Range("E48").Value = "'13.83%"
$ws.Range("G48").Value = "'19"
$ws.Range("D49").Value = "'0.001998"
$ws.Range("E49").Value = "'-30.98%"
$ws.Range("G49").Value = "'19"
$ws.Range("D50").Value = "'0.00002098"
$ws.Range("E50").Value = "'0.03%"
$ws.Range("G50").Value = "'19"
$ws.Range("D51").Value = "'0.0001998"
$ws.Range("E51").Value = "'0.03%"
$ws.Range("G51").Value = "'19"
